$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sign-in/sign-out log rows for Tuesday 2023-10-31
$newRows = @(
    @("Tuesday", "2023-10-31", "jesse", "Jesse Febian", "Safety", "Safety Officer", "17:54:37", "Signed-In"),
    @("Tuesday", "2023-10-31", "headt", "Head Teacher", "Admin", "Head Teacher", "18:02:14", "Signed-In"),
    @("Tuesday", "2023-10-31", "headt", "Head Teacher", "Admin", "Head Teacher", "18:05:50", "Signed-Out"),
    @("Tuesday", "2023-10-31", "jesse", "Jesse Febian", "Safety", "Safety Officer", "18:06:12", "Signed-Out")
)

$startRow = 7
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $col = $c + 1
        $text = $values[$c]
        $cell = $ws.Cells.Item($row, $col)
        if ($col -eq 2) {
            # Column B holds a date-looking string ("2023-10-31"). Assigning
            # it directly via .Value lets Excel auto-convert it to a date
            # serial number (and stamps a number-format style on the cell).
            # Route it through a text formula, then paste-special as values,
            # so it lands as plain text (shared string) with no style change.
            $cell.Formula = '="' + $text + '"'
            $cell.Copy()
            $cell.PasteSpecial(-4163)
        } else {
            $cell.Value = $text
        }
    }
}
